$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 233.33333
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 200
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -30
$ws.Range("N12").Value = -640
$ws.Range("H17").Value = 1099.0864
$ws.Range("I17").Value = 977.625
$ws.Range("J17").Value = 1150.228
$ws.Range("K17").Value = 2932.875
$ws.Range("L17").Value = 3450.684
$ws.Range("M17").Value = -2764.875
$ws.Range("N17").Value = -3786.684
$ws.Range("H33").Value = 876
$ws.Range("I33").Value = 1330.1666
$ws.Range("J33").Value = 194.75
$ws.Range("K33").Value = 1330.1666
$ws.Range("L33").Value = 194.75
$ws.Range("M33").Value = -1101.1666
$ws.Range("N33").Value = -652.75
$ws.Range("H43").Value = 1877.2307
$ws.Range("I43").Value = 1475
$ws.Range("J43").Value = 1997.9
$ws.Range("K43").Value = 1475
$ws.Range("L43").Value = 1997.9
$ws.Range("M43").Value = -1406
$ws.Range("N43").Value = -2135.9
$ws.Range("H55").Value = 62500920
$ws.Range("I55").Value = 83334400
$ws.Range("J55").Value = 493
$ws.Range("K55").Value = 83334400
$ws.Range("L55").Value = 493
$ws.Range("M55").Value = -83334186
$ws.Range("N55").Value = -921
$ws.Range("H100").Value = 3190.9092
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 3637.5
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 3637.5
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -4719.5
$ws.Range("H112").Value = 3925.4375
$ws.Range("I112").Value = 980
$ws.Range("J112").Value = 3988.1064
$ws.Range("K112").Value = 2940
$ws.Range("L112").Value = 11964.3192
$ws.Range("M112").Value = -1832
$ws.Range("N112").Value = -14180.3192
$ws.Range("H132").Value = 2952.3784
$ws.Range("I132").Value = 2795.121
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 8385.363000000001
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -5855.363000000001
$ws.Range("N132").Value = -17809.25
$ws.Range("H134").Value = 75156
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 75156
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 75156
$ws.Range("N134").Value = -85296

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 2400
$ws.Range("I22").Value = 1800
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -1501
$ws.Range("N22").Value = -3598
$ws.Range("H32").Value = 4099.76
$ws.Range("I32").Value = 4099.76
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4099.76
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3812.76
$ws.Range("H45").Value = 1874
$ws.Range("I45").Value = 968
$ws.Range("J45").Value = 2931
$ws.Range("K45").Value = 968
$ws.Range("L45").Value = 2931
$ws.Range("M45").Value = -591
$ws.Range("N45").Value = -3685
$ws.Range("H109").Value = 93750
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 93750
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 93750
$ws.Range("N109").Value = -96524
$ws.Range("H132").Value = 1186260.1
$ws.Range("I132").Value = 2027.2449
$ws.Range("J132").Value = 4812973
$ws.Range("K132").Value = 6081.7347
$ws.Range("L132").Value = 14438919
$ws.Range("M132").Value = -3551.7347
$ws.Range("N132").Value = -14443979

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29413286
$ws.Range("I20").Value = 1484.5
$ws.Range("J20").Value = 100001610
$ws.Range("K20").Value = 1484.5
$ws.Range("L20").Value = 100001610
$ws.Range("M20").Value = -1237.5
$ws.Range("N20").Value = -100002104

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6011.7764
$ws.Range("I31").Value = 2814.16
$ws.Range("J31").Value = 7344.1167
$ws.Range("K31").Value = 2814.16
$ws.Range("L31").Value = 7344.1167
$ws.Range("M31").Value = -2519.16
$ws.Range("N31").Value = -7934.1167
$ws.Range("H34").Value = 6011.7764
$ws.Range("I34").Value = 2814.16
$ws.Range("J34").Value = 7344.1167
$ws.Range("K34").Value = 2814.16
$ws.Range("L34").Value = 7344.1167
$ws.Range("M34").Value = -2612.16
$ws.Range("N34").Value = -7748.1167
$ws.Range("H98").Value = 37695
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 37695
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 37695
$ws.Range("N98").Value = -42187
$ws.Range("H134").Value = 7148192.5
$ws.Range("I134").Value = 10422646
$ws.Range("J134").Value = 3930.2727
$ws.Range("K134").Value = 31267938
$ws.Range("L134").Value = 11790.8181
$ws.Range("M134").Value = -31265403
$ws.Range("N134").Value = -16860.8181

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1770.8823
$ws.Range("I55").Value = 1400
$ws.Range("J55").Value = 1794.0625
$ws.Range("K55").Value = 4200
$ws.Range("L55").Value = 5382.1875
$ws.Range("M55").Value = -4023
$ws.Range("N55").Value = -5736.1875
$ws.Range("H117").Value = 2595.9614
$ws.Range("I117").Value = 1792.4445
$ws.Range("J117").Value = 3021.353
$ws.Range("K117").Value = 5377.333500000001
$ws.Range("L117").Value = 9064.059000000001
$ws.Range("M117").Value = -1935.333500000001
$ws.Range("N117").Value = -15948.059

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 423.6
$ws.Range("I107").Value = 305
$ws.Range("J107").Value = 749.75
$ws.Range("K107").Value = 305
$ws.Range("L107").Value = 749.75
$ws.Range("M107").Value = 1615
$ws.Range("N107").Value = -4589.75
$ws.Range("H122").Value = 3185.7144
$ws.Range("I122").Value = 4025
$ws.Range("J122").Value = 2066.6667
$ws.Range("K122").Value = 12075
$ws.Range("L122").Value = 6200.000100000001
$ws.Range("M122").Value = -9625
$ws.Range("N122").Value = -11100.0001
$ws.Range("H123").Value = 8565.237999999999
$ws.Range("I123").Value = 3000
$ws.Range("J123").Value = 9874.706
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 9874.706
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -14774.706

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 100014
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 100014
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 100014
$ws.Range("N43").Value = -100400
$ws.Range("H61").Value = 5513.125
$ws.Range("I61").Value = 2866.6667
$ws.Range("J61").Value = 7101
$ws.Range("K61").Value = 2866.6667
$ws.Range("L61").Value = 7101
$ws.Range("M61").Value = -2664.6667
$ws.Range("N61").Value = -7505
$ws.Range("H64").Value = 98150
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 98150
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 98150
$ws.Range("N64").Value = -98600
$ws.Range("H67").Value = 98150
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 98150
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 98150
$ws.Range("N67").Value = -99710
$ws.Range("H76").Value = 62262.668
$ws.Range("I76").Value = 22500
$ws.Range("J76").Value = 82144
$ws.Range("K76").Value = 22500
$ws.Range("L76").Value = 82144
$ws.Range("M76").Value = -22162
$ws.Range("N76").Value = -82820
$ws.Range("H79").Value = 62262.668
$ws.Range("I79").Value = 22500
$ws.Range("J79").Value = 82144
$ws.Range("K79").Value = 22500
$ws.Range("L79").Value = 82144
$ws.Range("M79").Value = -21330
$ws.Range("N79").Value = -84484
$ws.Range("H113").Value = 5513.125
$ws.Range("I113").Value = 2866.6667
$ws.Range("J113").Value = 7101
$ws.Range("K113").Value = 2866.6667
$ws.Range("L113").Value = 7101
$ws.Range("M113").Value = -696.6667000000002
$ws.Range("N113").Value = -11441

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 500005000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 500005000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 500005000
$ws.Range("N32").Value = -500005634
$ws.Range("H123").Value = 98414.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 98414.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 98414.5
$ws.Range("N123").Value = -108214.5

Write-Host "Applied all profit updates"